$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item(2)

# Update rows 2-4 with new values
$ws.Range("A2").Value = 29
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 2200
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1

$ws.Range("A3").Value = 30
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 5998
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 2

$ws.Range("A4").Value = 31
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 10000
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 10

# Remove rows 5 through 10 (old extra data) entirely
$ws.Rows("5:10").Delete()
